$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.634423333333332
$ws.Range("H2").Value = 25.90327
$ws.Range("I2").Value = 0.8921355692341975
$ws.Range("J2").Value = 0.9242739139243521
$ws.Range("M2").Value = 36.48539666666667
$ws.Range("N2").Value = 109.45619
$ws.Range("O2").Value = 0.4260639713374229
$ws.Range("P2").Value = 0.4324607845540777
$ws.Range("Q2").Value = 315.0303603045888
$ws.Range("R2").Value = 2835.2732427413
$ws.Range("S2").Value = 0.3801068235992945
$ws.Range("T2").Value = 0.3997122219585934
$ws.Range("G3").Value = 8.634423333333332
$ws.Range("H3").Value = 25.90327
$ws.Range("I3").Value = 0.8921355692341975
$ws.Range("J3").Value = 0.9242739139243521
$ws.Range("O3").Value = 0.1743777127077069
$ws.Range("P3").Value = 0.1769957741547643
$ws.Range("Q3").Value = 128.9343323045111
$ws.Range("R3").Value = 1160.4089907406
$ws.Range("S3").Value = 0.1555685599882474
$ws.Range("T3").Value = 0.1635925769260947
$ws.Range("G4").Value = 8.634423333333332
$ws.Range("H4").Value = 25.90327
$ws.Range("I4").Value = 0.8921355692341975
$ws.Range("J4").Value = 0.9242739139243521
$ws.Range("M4").Value = 9.680823666666667
$ws.Range("N4").Value = 29.042471
$ws.Range("O4").Value = 0.1130493445068016
$ws.Range("P4").Value = 0.1147466378470605
$ws.Range("Q4").Value = 83.58832975335221
$ws.Range("R4").Value = 752.2949677801699
$ws.Range("S4").Value = 0.1008553413131284
$ws.Range("T4").Value = 0.1060573240725628
$ws.Range("G5").Value = 8.634423333333332
$ws.Range("H5").Value = 25.90327
$ws.Range("I5").Value = 0.8921355692341975
$ws.Range("J5").Value = 0.9242739139243521
$ws.Range("M5").Value = 3.79999
$ws.Range("N5").Value = 7.59998
$ws.Range("O5").Value = 0.04437498227672168
$ws.Range("P5").Value = 0.0300274777826206
$ws.Range("Q5").Value = 32.81072232243333
$ws.Range("R5").Value = 196.8643339346
$ws.Range("S5").Value = 0.03958850007320052
$ws.Range("T5").Value = 0.02775361441541927
$ws.Range("G6").Value = 8.634423333333332
$ws.Range("H6").Value = 25.90327
$ws.Range("I6").Value = 0.8921355692341975
$ws.Range("J6").Value = 0.9242739139243521
$ws.Range("M6").Value = 20.734808
$ws.Range("N6").Value = 62.204424
$ws.Range("O6").Value = 0.242133989171347
$ws.Range("P6").Value = 0.245769325661477
$ws.Range("Q6").Value = 179.0331100073867
$ws.Range("R6").Value = 1611.29799006648
$ws.Range("S6").Value = 0.2160163442603267
$ws.Range("T6").Value = 0.227158176551682
$ws.Range("I7").Value = 0.003550079123985249
$ws.Range("J7").Value = 0.003677967385028239
$ws.Range("M7").Value = 36.48539666666667
$ws.Range("N7").Value = 109.45619
$ws.Range("O7").Value = 0.4260639713374229
$ws.Range("P7").Value = 0.4324607845540777
$ws.Range("Q7").Value = 1.25360174407
$ws.Range("R7").Value = 11.28241569663
$ws.Range("S7").Value = 0.001512560810127234
$ws.Range("T7").Value = 0.001590576660893622
$ws.Range("I8").Value = 0.003550079123985249
$ws.Range("J8").Value = 0.003677967385028239
$ws.Range("O8").Value = 0.1743777127077069
$ws.Range("P8").Value = 0.1769957741547643
$ws.Range("S8").Value = 0.0006190546775719273
$ws.Range("T8").Value = 0.0006509846846290473
$ws.Range("I9").Value = 0.003550079123985249
$ws.Range("J9").Value = 0.003677967385028239
$ws.Range("M9").Value = 9.680823666666667
$ws.Range("N9").Value = 29.042471
$ws.Range("O9").Value = 0.1130493445068016
$ws.Range("P9").Value = 0.1147466378470605
$ws.Range("Q9").Value = 0.332623420363
$ws.Range("R9").Value = 2.993610783267
$ws.Range("S9").Value = 0.0004013341179138129
$ws.Range("T9").Value = 0.0004220343915431356
$ws.Range("I10").Value = 0.003550079123985249
$ws.Range("J10").Value = 0.003677967385028239
$ws.Range("M10").Value = 3.79999
$ws.Range("N10").Value = 7.59998
$ws.Range("O10").Value = 0.04437498227672168
$ws.Range("P10").Value = 0.0300274777826206
$ws.Range("Q10").Value = 0.13056385641
$ws.Range("R10").Value = 0.7833831384600001
$ws.Range("S10").Value = 0.000157534698207805
$ws.Range("T10").Value = 0.0001104400839391386
$ws.Range("I11").Value = 0.003550079123985249
$ws.Range("J11").Value = 0.003677967385028239
$ws.Range("M11").Value = 20.734808
$ws.Range("N11").Value = 62.204424
$ws.Range("O11").Value = 0.242133989171347
$ws.Range("P11").Value = 0.245769325661477
$ws.Range("Q11").Value = 0.7124272680720001
$ws.Range("R11").Value = 6.411845412648001
$ws.Range("S11").Value = 0.0008595948201644694
$ws.Range("T11").Value = 0.0009039315640232961
$ws.Range("G12").Value = 1.0095935
$ws.Range("H12").Value = 2.019187
$ws.Range("I12").Value = 0.1043143516418173
$ws.Range("J12").Value = 0.07204811869061979
$ws.Range("M12").Value = 36.48539666666667
$ws.Range("N12").Value = 109.45619
$ws.Range("O12").Value = 0.4260639713374229
$ws.Range("P12").Value = 0.4324607845540777
$ws.Range("Q12").Value = 36.83541931958833
$ws.Range("R12").Value = 221.01251591753
$ws.Range("S12").Value = 0.04444458692800111
$ws.Range("T12").Value = 0.03115798593459074
$ws.Range("G13").Value = 1.0095935
$ws.Range("H13").Value = 2.019187
$ws.Range("I13").Value = 0.1043143516418173
$ws.Range("J13").Value = 0.07204811869061979
$ws.Range("O13").Value = 0.1743777127077069
$ws.Range("P13").Value = 0.1769957741547643
$ws.Range("Q13").Value = 15.07584916747667
$ws.Range("R13").Value = 90.45509500486001
$ws.Range("S13").Value = 0.01819009804188753
$ws.Range("T13").Value = 0.01275221254404059
$ws.Range("G14").Value = 1.0095935
$ws.Range("H14").Value = 2.019187
$ws.Range("I14").Value = 0.1043143516418173
$ws.Range("J14").Value = 0.07204811869061979
$ws.Range("M14").Value = 9.680823666666667
$ws.Range("N14").Value = 29.042471
$ws.Range("O14").Value = 0.1130493445068016
$ws.Range("P14").Value = 0.1147466378470605
$ws.Range("Q14").Value = 9.773696648512834
$ws.Range("R14").Value = 58.642179891077
$ws.Range("S14").Value = 0.01179266907575946
$ws.Range("T14").Value = 0.008267279382954582
$ws.Range("G15").Value = 1.0095935
$ws.Range("H15").Value = 2.019187
$ws.Range("I15").Value = 0.1043143516418173
$ws.Range("J15").Value = 0.07204811869061979
$ws.Range("M15").Value = 3.79999
$ws.Range("N15").Value = 7.59998
$ws.Range("O15").Value = 0.04437498227672168
$ws.Range("P15").Value = 0.0300274777826206
$ws.Range("Q15").Value = 3.836445204065
$ws.Range("R15").Value = 15.34578081626
$ws.Range("S15").Value = 0.004628947505313358
$ws.Range("T15").Value = 0.002163423283262198
$ws.Range("G16").Value = 1.0095935
$ws.Range("H16").Value = 2.019187
$ws.Range("I16").Value = 0.1043143516418173
$ws.Range("J16").Value = 0.07204811869061979
$ws.Range("M16").Value = 20.734808
$ws.Range("N16").Value = 62.204424
$ws.Range("O16").Value = 0.242133989171347
$ws.Range("P16").Value = 0.245769325661477
$ws.Range("Q16").Value = 20.933727380548
$ws.Range("R16").Value = 125.602364283288
$ws.Range("S16").Value = 0.02525805009085588
$ws.Range("T16").Value = 0.01770721754577168
